$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-add a previously removed row of temporal oil import data at the bottom
# of the table. Inserting (rather than just writing into the blank row)
# pulls down the formatting of row 83 automatically, so the new cells pick
# up the same cell styles (wrap-text / font flags) as their neighbours.
$ws.Rows.Item(84).Insert(-4121)

$ws.Range("A84").Value = "oil_imports"
$ws.Range("B84").Value = "csv"
$ws.Range("C84").Value = "Socioeconomic Data"
$ws.Range("D84").Value = "monthly oil and gas imports"
$ws.Range("E84").Value = "na"
$ws.Range("F84").Value = "Indonesia"
$ws.Range("G84").Value = "https://drive.google.com/drive/folders/1gg1jPYMPD0pWS5mMMmnUMTnYY39pWnFI"

# Column H isn't part of this row (unlike its neighbours) - strip the
# inherited format/content from the insert so the cell stays empty.
$ws.Range("A84").Copy($ws.Range("H84"))
$ws.Range("H84").ClearContents()

# Update the saved selection to mirror the view captured on save. (The
# header row stays frozen - ySplit is untouched - only where the cursor
# was left changes.)
$ws.Activate()
$ws.Range("G92").Select()
